$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A13").Value = "Dowel alignment pins for gear"
$ws.Range("D13").Value = "http://www.mcmaster.com/#97155a426/=xwsf2k"
$ws.Range("C13").Value = "1 pack"

$ws.Range("E5").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Value = 3

$ws.Range("B14").Select()
